# Daily data entry: append the two new station readings for 2025-08-29
# (45898) and 2025-08-30 (45899), mirroring the existing "date / station"
# pattern already used in column A/B of the sheet. The service-fee /
# revenue / order-count columns (C:D:E:F) for these two new days are not
# yet known, so they are intentionally left blank (the row only gets
# A/B filled in, same as every other "pending" row at the bottom of the
# sheet before this commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A58").Value = 45898
$ws.Range("B58").Value = "四方坪站"

$ws.Range("A59").Value = 45898
$ws.Range("B59").Value = "高岭站"

$ws.Range("A60").Value = 45899
$ws.Range("B60").Value = "四方坪站"

$ws.Range("A61").Value = 45899
$ws.Range("B61").Value = "高岭站"

# Move the selection down to where the new data was entered and scroll
# the window so that area is in view, matching the author's on-screen
# position after typing in the new rows.
$ws.Range("D62").Select()
